$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ID" column (old column A) was removed from the table; Arc (old B)
# and Description (old C) shift left to become A and B.
$ws.Columns.Item(1).Delete() | Out-Null

# Row 5 (Marineford) gets a slightly smaller row height after the reflow.
$ws.Rows.Item(5).RowHeight = 152.4

# Update the view: zoomed to 69% and selection moved to R3 (no more
# topLeftCell freeze on A3).
$excel.ActiveWindow.Zoom = 69
$ws.Range("R3").Select() | Out-Null
